$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column B per the commit:
#  Account Row (B4): 6 -> 3
#  Wealth Row (B5): 7 -> 3
#  Records Row (B6): 3 -> 2
#  Records Banks Column (B7): "J" -> "G"
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = "G"
